$d = $word.ActiveDocument

# --- Update the "Body Text" style: add 1.5x-ish (360 twips / auto) line
# spacing on top of the existing before/after spacing, and set the run
# size to 11pt (22 half-points). ---
$bodyText = $d.Styles("Body Text")
$bodyText.ParagraphFormat.LineSpacingRule = 5   # wdLineSpaceAuto-with-multiple -> w:lineRule="auto"
$bodyText.ParagraphFormat.LineSpacing = 18      # 18pt = 360 twips -> w:line="360"
$bodyText.Font.Size = 11                        # -> w:sz w:val="22"

# --- Update the linked "Body Text Char" character style to match: add the
# same 11pt run size, keeping its existing Times New Roman font. ---
$bodyTextChar = $d.Styles("Body Text Char")
$bodyTextChar.Font.Size = 11

# --- Add the new "Body Text 2" paragraph style (based on Normal, linked to
# a new "Body Text 2 Char" character style). ---
$bodyText2 = $d.Styles.Add("Body Text 2", 1)     # wdStyleTypeParagraph
$bodyText2.BaseStyle = "Normal"
$bodyText2.LinkStyle = "BodyText2Char"
$bodyText2.ParagraphFormat.SpaceAfter = 6        # 6pt = 120 twips -> w:after="120"
$bodyText2.ParagraphFormat.LineSpacingRule = 5
$bodyText2.ParagraphFormat.LineSpacing = 18      # -> w:line="360" w:lineRule="auto"

# --- Add the new "Body Text 2 Char" character style linked back to it. ---
$bodyText2Char = $d.Styles.Add("Body Text 2 Char", 2)   # wdStyleTypeCharacter
$bodyText2Char.BaseStyle = "DefaultParagraphFont"
$bodyText2Char.LinkStyle = "BodyText2"
$bodyText2Char.Font.Name = "Times New Roman"
